$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.025641
$ws.Range("D2").Value = 0.101788
$ws.Range("E2").Value = 0.749993

$ws.Range("B3").Value = 1505.652417
$ws.Range("D3").Value = 4.781625
$ws.Range("E3").Value = 0.009261

$ws.Range("B4").Value = 35109.453286
$ws.Range("C4").Value = 223

$ws.Range("G5").Value = 4.732538
$ws.Range("H5").Value = -0.747439
$ws.Range("I5").Value = 10.212514
$ws.Range("J5").Value = 0.105662

$ws.Range("G6").Value = -0.641026
$ws.Range("H6").Value = -6.5418
$ws.Range("I6").Value = 5.259749
$ws.Range("J6").Value = 0.964437

$ws.Range("G7").Value = -5.373563
$ws.Range("H7").Value = -9.83461
$ws.Range("I7").Value = -0.912516
$ws.Range("J7").Value = 0.013528
